# Weekly fruit/vegetable price update: insert a new observation row into the
# "Femacal de La Calera - Zapallo italiano" data table.
#
# The new record is inserted as row 143 (pushing the previous rows 143..228
# down to 144..229), growing the used range from A1:R228 to A1:R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 143; this shifts existing rows 143-228
# down to 144-229 and extends the sheet dimension automatically.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record's data.
$ws.Range("A143").Value = 3
$ws.Range("B143").Value = "Femacal de La Calera"
$ws.Range("C143").Value = "Coquimbo"
$ws.Range("D143").Value = 44488
$ws.Range("E143").Value = 5
$ws.Range("F143").Value = 100112032
$ws.Range("G143").Value = "Zapallo italiano"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 200
$ws.Range("K143").Value = 9000
$ws.Range("L143").Value = 10000
$ws.Range("M143").Value = 9525
$ws.Range("N143").Value = "`$/caja 70 unidades"
$ws.Range("O143").Value = "Región de Arica y Parinacota"
$ws.Range("P143").Value = 136
$ws.Range("Q143").Value = 70
$ws.Range("R143").Value = "Hortaliza"
